$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Electrical Equipment"
$ws.Range("B2").Value = "'1"

$ws.Range("A3").Value = "Labor — Contract"
$ws.Range("B3").Value = "'1"

$ws.Range("A4").Value = "Labor — Services"
$ws.Range("B4").Value = "'1"

$ws.Range("A5").Value = "Oil  Related Products"
$ws.Range("B5").Value = "'1"

$ws.Range("A6").Value = "Fuel"
$ws.Range("B6").Value = "'2"

$ws.Range("A7").Value = "Diesel Fuel"
$ws.Range("B7").Value = "'3"

$ws.Range("A8").Value = "Electrical Components"
$ws.Range("B8").Value = "'32"

$ws.Range("A9").Value = "Labor"
$ws.Range("B9").Value = "'34"

$ws.Range("A10").Value = "Lumber"
$ws.Range("B10").Value = "'4"

$ws.Range("A11").Value = "Gasoline"
$ws.Range("B11").Value = "'8"

$ws.Range("A12").Value = " Polyvinyl Chloride"
$ws.Range("B12").Value = "PVC) Products"

$ws.Range("B2:B11").Style = "Normal"
